# The deck's two theme parts (theme1.xml, used by the slide master / overall
# design, and theme2.xml, used only by the notes master) had their palettes
# swapped: theme1 switched from the "Integral" (Red Violet) palette to the
# stock "Office Theme" (Office) palette, and theme2 switched the other way.
#
# Font scheme and format scheme (fills/lines/effects) are identical between
# the two theme parts, so only the 12-colour scheme actually changes.
#
# Apply the new "Office" palette to the presentation's theme colour scheme
# (maps onto ppt/theme/theme1.xml's <a:clrScheme>) via ThemeColorScheme,
# which is the documented way to edit theme colours through this object
# model: ThemeColorScheme.Colors(i).RGB.
#
# RGB() packs as 0xBBGGRR (red + green*256 + blue*65536), matching the
# VBA RGB() helper.
function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Office theme colour scheme, in DrawingML clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
